# Natmi following Dr Hou advice
# Rebuild the LR-pairs table for Vwf-Tnfrsf11b to include the FAPs cluster
# as both a sending and a target cluster, alongside the existing ECs and sCs
# clusters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> Vwf/Tnfrsf11b -> FAPs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Vwf"
$ws.Range("C2").Value = "Tnfrsf11b"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 58.356725
$ws.Range("H2").Value = 175.070175
$ws.Range("I2").Value = 0.9850455027887332
$ws.Range("J2").Value = 0.9850455027887334
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.095195666666667
$ws.Range("N2").Value = 6.285587
$ws.Range("O2").Value = 0.8546922300706357
$ws.Range("P2").Value = 0.8546922300706358
$ws.Range("Q2").Value = 122.2687573408583
$ws.Range("R2").Value = 1100.418816067725
$ws.Range("S2").Value = 0.8419107374995529
$ws.Range("T2").Value = 0.8419107374995531

# Row 3: ECs -> Vwf/Tnfrsf11b -> sCs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Vwf"
$ws.Range("C3").Value = "Tnfrsf11b"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 58.356725
$ws.Range("H3").Value = 175.070175
$ws.Range("I3").Value = 0.9850455027887332
$ws.Range("J3").Value = 0.9850455027887334
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.356208
$ws.Range("N3").Value = 1.068624
$ws.Range("O3").Value = 0.1453077699293643
$ws.Range("P3").Value = 0.1453077699293643
$ws.Range("Q3").Value = 20.7871322988
$ws.Range("R3").Value = 187.0841906892
$ws.Range("S3").Value = 0.1431347652891802
$ws.Range("T3").Value = 0.1431347652891802

# Row 4: FAPs -> Vwf/Tnfrsf11b -> FAPs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Vwf"
$ws.Range("C4").Value = "Tnfrsf11b"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.3120586666666667
$ws.Range("H4").Value = 0.936176
$ws.Range("I4").Value = 0.005267464653066949
$ws.Range("J4").Value = 0.00526746465306695
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.095195666666667
$ws.Range("N4").Value = 6.285587
$ws.Range("O4").Value = 0.8546922300706357
$ws.Range("P4").Value = 0.8546922300706358
$ws.Range("Q4").Value = 0.6538239661457778
$ws.Range("R4").Value = 5.884415695312
$ws.Range("S4").Value = 0.004502061111148038
$ws.Range("T4").Value = 0.004502061111148039

# Row 5: FAPs -> Vwf/Tnfrsf11b -> sCs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Vwf"
$ws.Range("C5").Value = "Tnfrsf11b"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.3120586666666667
$ws.Range("H5").Value = 0.936176
$ws.Range("I5").Value = 0.005267464653066949
$ws.Range("J5").Value = 0.00526746465306695
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.356208
$ws.Range("N5").Value = 1.068624
$ws.Range("O5").Value = 0.1453077699293643
$ws.Range("P5").Value = 0.1453077699293643
$ws.Range("Q5").Value = 0.111157793536
$ws.Range("R5").Value = 1.000420141824
$ws.Range("S5").Value = 0.0007654035419189109
$ws.Range("T5").Value = 0.000765403541918911

# Row 6: sCs -> Vwf/Tnfrsf11b -> FAPs
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Vwf"
$ws.Range("C6").Value = "Tnfrsf11b"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.5738856666666666
$ws.Range("H6").Value = 1.721657
$ws.Range("I6").Value = 0.00968703255819983
$ws.Range("J6").Value = 0.009687032558199832
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.095195666666667
$ws.Range("N6").Value = 6.285587
$ws.Range("O6").Value = 0.8546922300706357
$ws.Range("P6").Value = 0.8546922300706358
$ws.Range("Q6").Value = 1.202402761962111
$ws.Range("R6").Value = 10.821624857659
$ws.Range("S6").Value = 0.008279431459934667
$ws.Range("T6").Value = 0.008279431459934671

# Row 7: sCs -> Vwf/Tnfrsf11b -> sCs
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Vwf"
$ws.Range("C7").Value = "Tnfrsf11b"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.5738856666666666
$ws.Range("H7").Value = 1.721657
$ws.Range("I7").Value = 0.00968703255819983
$ws.Range("J7").Value = 0.009687032558199832
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.356208
$ws.Range("N7").Value = 1.068624
$ws.Range("O7").Value = 0.1453077699293643
$ws.Range("P7").Value = 0.1453077699293643
$ws.Range("Q7").Value = 0.204422665552
$ws.Range("R7").Value = 1.839803989968
$ws.Range("S7").Value = 0.001407601098265162
$ws.Range("T7").Value = 0.001407601098265162
